$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Tue Feb 11 20:22:07 EST 2025"
$ws.Range("B3").Value = "Tue Feb 11 20:22:21 EST 2025"
$ws.Range("B4").Value = "Tue Feb 11 20:22:35 EST 2025"
